# BIS-1002: removed "Internal Assignment" column from export.
#
# The "Internal Assignment" column (column O) is dropped from both
# property-table blocks on the sheet: the header cell (O4 / O12) that held
# the "Internal Assignment" label, and the data cells below it (O5:O7 /
# O13:O15) that held the per-property "FALSE" flag. Their styles stay in
# place; only the cell contents are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# First table (rows 4-7): clear the "Internal Assignment" header + values.
$ws.Range("O4:O7").ClearContents()

# Second table (rows 12-15): clear the "Internal Assignment" header + values.
$ws.Range("O12:O15").ClearContents()

# Leave the same selection behind that the authored change left in place.
$ws.Range("O4:O15").Select()
